$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new requirement row (row 8)
$ws.Range("A8").Value = "rml-io-r7"
$ws.Range("B8").Value = "Logical sources and logical targets may indicate relative paths to resources"
$ws.Range("C8").Value = "rml-io"

# Match the row height used by the other requirement rows
$ws.Rows.Item(8).RowHeight = 34

# Move the active selection to B9, as in the edited workbook
$ws.Range("B9").Select()
